$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Serping1"
$ws.Range("C2").Value = "Selp"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 29.82164133333333
$ws.Range("H2").Value = 89.464924
$ws.Range("I2").Value = 0.02335016309719764
$ws.Range("J2").Value = 0.02335016309719765
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 15.960008
$ws.Range("N2").Value = 47.880024
$ws.Range("O2").Value = 0.9899803616776065
$ws.Range("P2").Value = 0.9899803616776066
$ws.Range("Q2").Value = 475.9536342531306
$ws.Range("R2").Value = 4283.582708278176
$ws.Range("S2").Value = 0.02311620290819482
$ws.Range("T2").Value = 0.02311620290819483

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Serping1"
$ws.Range("C3").Value = "Selp"
$ws.Range("D3").Value = "M2"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 29.82164133333333
$ws.Range("H3").Value = 89.464924
$ws.Range("I3").Value = 0.02335016309719764
$ws.Range("J3").Value = 0.02335016309719765
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.06197766666666666
$ws.Range("N3").Value = 0.185933
$ws.Range("O3").Value = 0.003844401134548353
$ws.Range("P3").Value = 0.003844401134548354
$ws.Range("Q3").Value = 1.848275746010222
$ws.Range("R3").Value = 16.634481714092
$ws.Range("S3").Value = 0.000089767393502755699064010741
$ws.Range("T3").Value = 0.000089767393502755726169065054

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Serping1"
$ws.Range("C4").Value = "Selp"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 29.82164133333333
$ws.Range("H4").Value = 89.464924
$ws.Range("I4").Value = 0.02335016309719764
$ws.Range("J4").Value = 0.02335016309719765
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.09955433333333334
$ws.Range("N4").Value = 0.298663
$ws.Range("O4").Value = 0.006175237187845165
$ws.Range("P4").Value = 0.006175237187845166
$ws.Range("Q4").Value = 2.968873621845778
$ws.Range("R4").Value = 26.719862596612
$ws.Range("S4").Value = 0.0001441927955000647
$ws.Range("T4").Value = 0.0001441927955000648

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Serping1"
$ws.Range("C5").Value = "Selp"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1160.126729666667
$ws.Range("H5").Value = 3480.380189
$ws.Range("I5").Value = 0.9083721465342723
$ws.Range("J5").Value = 0.9083721465342726
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 15.960008
$ws.Range("N5").Value = 47.880024
$ws.Range("O5").Value = 0.9899803616776065
$ws.Range("P5").Value = 0.9899803616776066
$ws.Range("Q5").Value = 18515.63188649384
$ws.Range("R5").Value = 166640.6869784445
$ws.Range("S5").Value = 0.8992705861638627
$ws.Range("T5").Value = 0.899270586163863

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Serping1"
$ws.Range("C6").Value = "Selp"
$ws.Range("D6").Value = "M2"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1160.126729666667
$ws.Range("H6").Value = 3480.380189
$ws.Range("I6").Value = 0.9083721465342723
$ws.Range("J6").Value = 0.9083721465342726
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.06197766666666666
$ws.Range("N6").Value = 0.185933
$ws.Range("O6").Value = 0.003844401134548353
$ws.Range("P6").Value = 0.003844401134548354
$ws.Range("Q6").Value = 71.90194774237077
$ws.Range("R6").Value = 647.1175296813369
$ws.Range("S6").Value = 0.003492146910728479
$ws.Range("T6").Value = 0.003492146910728481

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Serping1"
$ws.Range("C7").Value = "Selp"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1160.126729666667
$ws.Range("H7").Value = 3480.380189
$ws.Range("I7").Value = 0.9083721465342723
$ws.Range("J7").Value = 0.9083721465342726
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.09955433333333334
$ws.Range("N7").Value = 0.298663
$ws.Range("O7").Value = 0.006175237187845165
$ws.Range("P7").Value = 0.006175237187845166
$ws.Range("Q7").Value = 115.4956431541452
$ws.Range("R7").Value = 1039.460788387307
$ws.Range("S7").Value = 0.005609413459681176
$ws.Range("T7").Value = 0.005609413459681178

# Row 8
$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "Serping1"
$ws.Range("C8").Value = "Selp"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.2607933333333334
$ws.Range("H8").Value = 0.7823800000000001
$ws.Range("I8").Value = 0.0002041995878070102
$ws.Range("J8").Value = 0.0002041995878070102
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 15.960008
$ws.Range("N8").Value = 47.880024
$ws.Range("O8").Value = 0.9899803616776065
$ws.Range("P8").Value = 0.9899803616776066
$ws.Range("Q8").Value = 4.162263686346668
$ws.Range("R8").Value = 37.46037317712
$ws.Range("S8").Value = 0.0002021535817916021
$ws.Range("T8").Value = 0.0002021535817916022

# Row 9
$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "Serping1"
$ws.Range("C9").Value = "Selp"
$ws.Range("D9").Value = "M2"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.2607933333333334
$ws.Range("H9").Value = 0.7823800000000001
$ws.Range("I9").Value = 0.0002041995878070102
$ws.Range("J9").Value = 0.0002041995878070102
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.06197766666666666
$ws.Range("N9").Value = 0.185933
$ws.Range("O9").Value = 0.003844401134548353
$ws.Range("P9").Value = 0.003844401134548354
$ws.Range("Q9").Value = 0.01616336228222222
$ws.Range("R9").Value = 0.14547026054
$ws.Range("S9").Value = 0.00000078502512703957611980165208
$ws.Range("T9").Value = 0.00000078502512703957622568077049

# Row 10
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Serping1"
$ws.Range("C10").Value = "Selp"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.2607933333333334
$ws.Range("H10").Value = 0.7823800000000001
$ws.Range("I10").Value = 0.0002041995878070102
$ws.Range("J10").Value = 0.0002041995878070102
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.09955433333333334
$ws.Range("N10").Value = 0.298663
$ws.Range("O10").Value = 0.006175237187845165
$ws.Range("P10").Value = 0.006175237187845166
$ws.Range("Q10").Value = 0.02596310643777779
$ws.Range("R10").Value = 0.23366795794
$ws.Range("S10").Value = 0.0000012609808883685029521747739
$ws.Range("T10").Value = 0.0000012609808883685040109659580

# Row 11
$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Serping1"
$ws.Range("C11").Value = "Selp"
$ws.Range("D11").Value = "ECs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 86.94000199999999
$ws.Range("H11").Value = 260.820006
$ws.Range("I11").Value = 0.0680734907807228
$ws.Range("J11").Value = 0.06807349078072282
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 15.960008
$ws.Range("N11").Value = 47.880024
$ws.Range("O11").Value = 0.9899803616776065
$ws.Range("P11").Value = 0.9899803616776066
$ws.Range("Q11").Value = 1387.563127440016
$ws.Range("R11").Value = 12488.06814696014
$ws.Range("S11").Value = 0.06739141902375717
$ws.Range("T11").Value = 0.0673914190237572

# Row 12
$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Serping1"
$ws.Range("C12").Value = "Selp"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 86.94000199999999
$ws.Range("H12").Value = 260.820006
$ws.Range("I12").Value = 0.0680734907807228
$ws.Range("J12").Value = 0.06807349078072282
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.06197766666666666
$ws.Range("N12").Value = 0.185933
$ws.Range("O12").Value = 0.003844401134548353
$ws.Range("P12").Value = 0.003844401134548354
$ws.Range("Q12").Value = 5.388338463955332
$ws.Range("R12").Value = 48.49504617559799
$ws.Range("S12").Value = 0.0002617018051900776
$ws.Range("T12").Value = 0.0002617018051900777

# Row 13
$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Serping1"
$ws.Range("C13").Value = "Selp"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 86.94000199999999
$ws.Range("H13").Value = 260.820006
$ws.Range("I13").Value = 0.0680734907807228
$ws.Range("J13").Value = 0.06807349078072282
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.09955433333333334
$ws.Range("N13").Value = 0.298663
$ws.Range("O13").Value = 0.006175237187845165
$ws.Range("P13").Value = 0.006175237187845166
$ws.Range("Q13").Value = 8.655253939108666
$ws.Range("R13").Value = 77.897285451978
$ws.Range("S13").Value = 0.0004203699517755544
$ws.Range("T13").Value = 0.0004203699517755546
